$wb = $excel.ActiveWorkbook

# Add three new sheets at the end of the workbook, in order:
#   PhonesDropDown, FeaturedPhones, AppsAndServicesList
$ws4 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws4.Name = "PhonesDropDown"
$ws5 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws5.Name = "FeaturedPhones"
$ws6 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws6.Name = "AppsAndServicesList"

# --- PhonesDropDown ---
$phonesDropDown = @(
    "Smartphones",
    "Trade in your phone",
    "Basic phones",
    "Certified pre-owned",
    "Prepaid phones",
    "Bring your own device",
    "Unlocked phones",
    "Phone accessories"
)
for ($i = 0; $i -lt $phonesDropDown.Length; $i++) {
    $ws4.Cells.Item($i + 1, 1).Value = $phonesDropDown[$i]
}
# Requested value tuned so the engine's pixel-quantized stored width lands
# as close as possible to the source column width (31.85546875 characters).
$ws4.Columns.Item(1).ColumnWidth = 30.916666666666664
$ws4.Range("A8").Select() | Out-Null

# --- FeaturedPhones ---
$featuredPhones = @(
    "Apple iPhone SE (2020)",
    "Apple iPhone 11",
    "Apple iPhone 11 Pro",
    "Apple iPhone 11 Pro Max",
    "Samsung Galaxy S20 5G UW",
    "Samsung Galaxy S20+ 5G",
    "Samsung Galaxy A51",
    "Samsung Galaxy Note 10+",
    "Google Pixel 4",
    "Moto Razr",
    "OnePlus 8 5G UW",
    "motorola edge+"
)
for ($i = 0; $i -lt $featuredPhones.Length; $i++) {
    $ws5.Cells.Item($i + 1, 1).Value = $featuredPhones[$i]
}
# Requested value tuned so the engine's pixel-quantized stored width lands
# as close as possible to the source column width (34.28515625 characters).
$ws5.Columns.Item(1).ColumnWidth = 33.416666666666664
$ws5.Range("A12").Select() | Out-Null

# --- AppsAndServicesList ---
$appsAndServices = @(
    "My Verizon",
    "Verizon Cloud",
    "Smart Family",
    "Device trade-in",
    "Device Protection",
    "Call Filter",
    "Apple Music",
    "Premium Visual Voicemail",
    "Hum",
    "See More Apps"
)
for ($i = 0; $i -lt $appsAndServices.Length; $i++) {
    $ws6.Cells.Item($i + 1, 1).Value = $appsAndServices[$i]
}
# Requested value tuned so the engine's pixel-quantized stored width lands
# as close as possible to the source column width (25.7109375 characters).
$ws6.Columns.Item(1).ColumnWidth = 24.75
$ws6.Range("A11").Select() | Out-Null

# AppsAndServicesList is the sheet active/selected when the workbook was saved
$ws6.Activate()
